$d = $word.ActiveDocument

# Step 1: find "Nombres" and replace with "nombre".
# Toggling Bold on/off around the Text assignment forces the engine to
# split this word into its own run (rather than merging it back into the
# surrounding text), matching how Word preserves run boundaries created
# by a discrete edit.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Nombres", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Text = "nombre"
$rng.Bold = 0

# Step 2: locate "BAJO" (immediately followed by the comma) in the same
# sentence and insert " separados por espacio" right after it, again
# toggling formatting to force a distinct run for the inserted text.
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("riesgo BAJO,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.MoveStart(1, 7)  # skip past "riesgo "
$rng2.MoveEnd(1, -1)   # drop the trailing comma from the match
$rng2.Collapse(0)      # collapse to the end (right after "BAJO")
$rng2.Bold = 1
$rng2.InsertBefore(" separados por espacio")
$rng2.Bold = 0
